$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$title = @'
낮은 타이어 공기압 경고가 발생한 경우의 조치
'@

$code = @'
HT203033
'@

$body = @'
낮은 타이어 공기압 경고가 발생한 경우의 조치
<br>
<br><img src="https://www.volvocars.com/images/support/imgdced27303b034a5cc0a801521024f995_1_--_--_VOICEpnghigh.png" height="80px" width="80px"
<br>시스템 경고등이 켜지고 타이어 공기압 부족이라는 메시지가 표시되면 타이어 공기압을 점검한 후 바로 잡아야 합니다.
<br>
<br>1. 차량의 시동을 끕니다.
<br>2. 타이어 공기압 게이지*로 타이어 네 개 모두의 타이어 공기압을 점검합니다.
<br>3. 타이어에 올바른 공기압으로 공기를 주입합니다. 공장 장착 타이어 권장 공기압이 표시된 운전석 도어 필라의 라벨을 참조합니다.
<br>4. 타이어 공기압을 조정한 후 항상 중앙 화면에서 새 타이어 공기압을 시스템에 저장하십시오. 이것은 차량이 켜져 있고, 정차한 상태 일 때만 실행할 수 있습니다.
<br>
<br>새 타이어 공기압 저장이 시작될 때까지 표시등이 꺼지지 않습니다.
<br>
<br>시스템에 새로운 기준치를 저장하려면 35km/h (22mph)를 넘는 속도로 차량을 몇 분간 주행 시킬 필요가 있습니다.
<br>
<br>참고
<br>타이어 공기압이 잘못 측정되는 것을 방지하기 위해 타이어 공기압은 상온 타이어 상태에서(타이어가 식었을 때) 점검하십시오. 상온 타이어는 온도가 대기 온도와 같은 타이어를 가리킵니다. 운전을 마친 후 차를 약 3시간 동안 세워 놓으면 타이어 온도가 대기 온도와 같아집니다. 차를 몇 킬로미터 운전하면 타이어 온도가 올라가 타이어 공기압이 높아지게 됩니다.
<br>
<br>참고
<br>• 타이어에 공기를 주입한 다음에는 노면에서 튀는 돌 등으로 밸브가 손상되지 않도록 밸브에 더스트 캡을 씌우십시오.
<br>• 더스트 캡은 플라스틱으로 된 것을 사용하십시오. 금속으로 된 것은 녹이 슬어 돌리기가 어렵게 될 수 있습니다.
<br>
<br><a style="color: red;">경고</a>
<br>• 타이어 공기압이 올바르지 않으면 타이어 고장이 발생하여 운전자가 차량에 대한 통제력을 상실할 수 있습니다.
<br>• 시스템은 갑작스러운 타이어 손상을 미리 알려줄 수 없습니다.
<br>
<br>* Indirect Tyre Pressure Monitoring System (iTPMS)
'@

$row = 45

$ws.Cells.Item($row, 1).Value = $title
$ws.Cells.Item($row, 2).Value = 45455
$ws.Cells.Item($row, 3).Value = $code
$ws.Cells.Item($row, 4).Value = $body
$ws.Cells.Item($row, 5).Value = "All"

$ws.Range("D45").WrapText = $true
$ws.Rows.Item($row).RowHeight = 409.6

$ws.Range("D45").Select() | Out-Null
